# Update crypto price/volume table with latest scraped values.
# GitHub Actions scheduled refresh - Fri Aug 25 21:08:38 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $range = $ws.Range($cellRef)
    # Force text storage so purely-numeric-looking strings (prices like
    # "218.13") keep their original inline-string type instead of being
    # auto-coerced to a number (which would also drop trailing zeros).
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

$ws.Range('D2').Value = '26.108.55'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.651.66'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  -0.31%  '
Set-TextCell 'D5' '218.13'
$ws.Range('E5').Value = '  +0.14%  '
Set-TextCell 'D6' '0.5292'
$ws.Range('E6').Value = '  +1.62%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('E8').Value = '  -2.12%  '
Set-TextCell 'D9' '0.06313'
$ws.Range('E9').Value = '  +0.28%  '
Set-TextCell 'D10' '20.36'
$ws.Range('E10').Value = '  -2.95%  '
Set-TextCell 'D11' '0.07753'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D12' '4.473'
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.652.42'
$ws.Range('E13').Value = '  -0.50%  '
Set-TextCell 'D14' '0.5456'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').Value = '0.0₅8124'
$ws.Range('E15').Value = '  -0.45%  '
Set-TextCell 'D16' '65.27'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '26.129.00'
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('E18').Value = '  -0.34%  '
Set-TextCell 'D19' '4.541'
$ws.Range('E19').Value = '  -2.35%  '
Set-TextCell 'D20' '193.78'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('E21').Value = '  -0.78%  '
Set-TextCell 'D22' '5.979'
$ws.Range('E22').Value = '  -1.41%  '
Set-TextCell 'D23' '1.004'
$ws.Range('E23').Value = '  -0.33%  '
Set-TextCell 'D24' '139.92'
$ws.Range('E24').Value = '  +1.14%  '
Set-TextCell 'D25' '0.1241'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  +0.82%  '
Set-TextCell 'D27' '16.15'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  +1.91%  '
Set-TextCell 'D29' '0.05912'
$ws.Range('E29').Value = '  -1.03%  '
Set-TextCell 'D30' '1.281'
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('E31').Value = '  -2.10%  '
Set-TextCell 'D32' '3.234'
$ws.Range('E32').Value = '  -2.27%  '
Set-TextCell 'D33' '1.546'
$ws.Range('E33').Value = '  -5.34%  '
Set-TextCell 'D34' '2.413'
$ws.Range('E34').Value = '  +0.08%  '
Set-TextCell 'D35' '0.9440'
$ws.Range('E35').Value = '  -3.30%  '
Set-TextCell 'D36' '2.759'
$ws.Range('E36').Value = '  -0.68%  '
Set-TextCell 'D37' '0.5640'
$ws.Range('E37').Value = '  -3.87%  '
Set-TextCell 'D38' '0.01603'
$ws.Range('E38').Value = '  +1.18%  '
Set-TextCell 'D39' '5.855'
$ws.Range('E39').Value = '  -1.25%  '
Set-TextCell 'D40' '0.8456'
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').Value = '1.009.31'
$ws.Range('E42').Value = '  -2.33%  '
Set-TextCell 'D43' '100.77'
$ws.Range('E43').Value = '  +1.34%  '
$ws.Range('D44').Value = '1.798.48'
$ws.Range('E44').Value = '  -0.09%  '
Set-TextCell 'D45' '56.81'
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('E46').Value = '  -2.60%  '
$ws.Range('E47').Value = '  +0.21%  '
Set-TextCell 'D48' '0.4289'
$ws.Range('E48').Value = '  +1.43%  '
Set-TextCell 'D49' '1.475'
$ws.Range('E49').Value = '  +1.54%  '
Set-TextCell 'D50' '0.05150'
Set-TextCell 'D51' '7.804'
$ws.Range('E51').Value = '  -3.39%  '
